$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Occupancy Rate" header that sat above the year row.
$ws.Range("B1").ClearContents()

# Remove the state-name labels in column A (the whole column's worth of data).
$ws.Range("A1:A54").ClearContents()

# Replace the "Avg." column header with a clearer label and re-flow the
# average column from a percentage display into a plain two-decimal number.
$ws.Range("K2").Value = "Occupancy Rate Avg"
$ws.Range("K2:K54").NumberFormat = "0.00"
$ws.Columns("K").ColumnWidth = 18

# The used range now starts at B2, so move the selection there.
$ws.Range("B1").Select() | Out-Null
